# Move and copy of if-then-else steps
# In the "DB" sheet, the "params" table description block (rows 9-14) gains
# a new explanatory column: the existing notes column (H) shifts right to
# column I, and a new note "parental_flag" is written into the freed-up
# column G on the params row (row 14). The same column shift (H -> I) also
# applies to the earlier "actions" table block (rows 5-11) for consistency.
#
# We use Range.Cut(destination) instead of Value-copy + Clear so that the
# underlying shared-string entries keep their original indices (a plain
# clear/overwrite would drop the no-longer-referenced string and silently
# renumber every shared string after it).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DB")

# Shift the "notes" column from H to I for the actions/steps/params blocks.
$ws.Range("H6").Cut($ws.Range("I6"))
$ws.Range("H10").Cut($ws.Range("I10"))
$ws.Range("H11").Cut($ws.Range("I11"))
$ws.Range("H14").Cut($ws.Range("I14"))

# New column/field note for the params table: "parental_flag".
$ws.Range("G14").Value = "parental_flag"
